$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style donor cells from the existing sheet (before any edits) ---
# A5/B5 currently hold the "highlighted" (Pendiente) style; A6 holds the plain date style.
$donorHighlightDate = $ws.Range("A5")
$donorHighlightText = $ws.Range("B5")
$donorPlainDate = $ws.Range("A6")

function Set-TextValue($range, [string]$text) {
    # Force a numeric-looking string to be stored as literal text (shared string),
    # without leaving behind any extra/unused cell style.
    $range.Formula = '="' + ($text -replace '"','""') + '"'
    $range.Copy()
    $range.PasteSpecial(-4163) | Out-Null
}

# Row 5
$ws.Range("A5").Value = 45732.44519675926
Set-TextValue $ws.Range("B5") "9000640"
Set-TextValue $ws.Range("C5") "475206"
$ws.Range("D5").Value = "Puerta3070 cerradura rota"
$ws.Range("E5").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A5").PasteSpecial(-4122) | Out-Null

# Row 6
$ws.Range("A6").Value = 45732.44519675926
Set-TextValue $ws.Range("B6") "8001903"
Set-TextValue $ws.Range("C6") "424102"
$ws.Range("D6").Value = "Sin leva 901"
$ws.Range("E6").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A6").PasteSpecial(-4122) | Out-Null

# Row 7
$ws.Range("A7").Value = 45732.44519675926
Set-TextValue $ws.Range("B7") "10000139"
Set-TextValue $ws.Range("C7") "480112"
$ws.Range("D7").Value = "Sin leva 3070"
$ws.Range("E7").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A7").PasteSpecial(-4122) | Out-Null

# Row 8
$ws.Range("A8").Value = 45732.44519675926
Set-TextValue $ws.Range("B8") "8002433"
Set-TextValue $ws.Range("C8") "107009"
$ws.Range("D8").Value = "Falta de leva de stacker 901"
$ws.Range("E8").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A8").PasteSpecial(-4122) | Out-Null

# Row 9
$ws.Range("A9").Value = 45732.44519675926
Set-TextValue $ws.Range("B9") "8001083"
Set-TextValue $ws.Range("C9") "459606"
$ws.Range("D9").Value = "Zona 58 falta leva 901"
$ws.Range("E9").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A9").PasteSpecial(-4122) | Out-Null

# Row 10
$ws.Range("A10").Value = 45732.44519675926
Set-TextValue $ws.Range("B10") "8002154"
Set-TextValue $ws.Range("C10") "479906"
$ws.Range("D10").Value = "Falta leva de stakers "
$ws.Range("E10").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A10").PasteSpecial(-4122) | Out-Null

# Row 11
$ws.Range("A11").Value = 45732.44519675926
Set-TextValue $ws.Range("B11") "8001435"
Set-TextValue $ws.Range("C11") "105502"
$ws.Range("D11").Value = "Puerta principal abre sola "
$ws.Range("E11").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A11").PasteSpecial(-4122) | Out-Null

# Row 12
$ws.Range("A12").Value = 45732.44519675926
Set-TextValue $ws.Range("B12") "8001851"
Set-TextValue $ws.Range("C12") "422704"
$ws.Range("D12").Value = "Revisar llave número 20 no abre"
$ws.Range("E12").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A12").PasteSpecial(-4122) | Out-Null

# Row 13
$ws.Range("A13").Value = 45732.44519675926
Set-TextValue $ws.Range("B13") "8001673"
Set-TextValue $ws.Range("C13") "308103"
$ws.Range("D13").Value = "Sin leva 901"
$ws.Range("E13").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A13").PasteSpecial(-4122) | Out-Null

# Row 14
$ws.Range("A14").Value = 45732.44519675926
Set-TextValue $ws.Range("B14") "11000260"
Set-TextValue $ws.Range("C14") "447605"
$ws.Range("D14").Value = "Sin leva de stacker "
$ws.Range("E14").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A14").PasteSpecial(-4122) | Out-Null

# Row 15
$ws.Range("A15").Value = 45732.44519675926
Set-TextValue $ws.Range("B15") "7001449"
Set-TextValue $ws.Range("C15") "302110"
$ws.Range("D15").Value = "Sin leva 901`n"
$ws.Range("E15").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A15").PasteSpecial(-4122) | Out-Null

# Row 16
$ws.Range("A16").Value = 45732.44519675926
Set-TextValue $ws.Range("B16") "8001671"
Set-TextValue $ws.Range("C16") "308101"
$ws.Range("D16").Value = "Sin leva 901"
$ws.Range("E16").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A16").PasteSpecial(-4122) | Out-Null

# Row 17
$ws.Range("A17").Value = 45732.44519675926
Set-TextValue $ws.Range("B17") "9001005"
Set-TextValue $ws.Range("C17") "105014"
$ws.Range("D17").Value = "Stacker sin leva "
$ws.Range("E17").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A17").PasteSpecial(-4122) | Out-Null

# Row 18
$ws.Range("A18").Value = 45732.44519675926
Set-TextValue $ws.Range("B18") "7000264"
Set-TextValue $ws.Range("C18") "426502"
$ws.Range("D18").Value = "Pta stacker rota se desarma"
$ws.Range("E18").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A18").PasteSpecial(-4122) | Out-Null

# Row 19
$ws.Range("A19").Value = 45732.44519675926
Set-TextValue $ws.Range("B19") "7000621"
Set-TextValue $ws.Range("C19") "116204"
$ws.Range("D19").Value = "Puerta principal"
$ws.Range("E19").Value = "Pendiente"
$donorHighlightDate.Copy()
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$donorHighlightText.Copy()
$ws.Range("B19:E19").PasteSpecial(-4122) | Out-Null

# Row 20
$ws.Range("A20").Value = 45732.44519675926
Set-TextValue $ws.Range("B20") "1000008"
Set-TextValue $ws.Range("C20") "463604"
$ws.Range("D20").Value = "Sin leva 901"
$ws.Range("E20").Value = "Extraída"
$donorPlainDate.Copy()
$ws.Range("A20").PasteSpecial(-4122) | Out-Null

# Row 15 contains an embedded newline; AutoFit keeps the row at the default height
# instead of leaving a stale explicit row-height override behind.
$ws.Rows(15).AutoFit() | Out-Null

$excel.CutCopyMode = 0